# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.878.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5015"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.48%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2561"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.259"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.642.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.864.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.02%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5407"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7843"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.51"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.902.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "197.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.369"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.898"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.958"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.869"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1137"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.813"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.237"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("E31").Value = "  -3.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.248"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.522"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.360"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8862"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.596"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.130.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5520"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01555"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.65%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.665"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.21%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.775.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4517"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05088"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
